$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-10 Sunday" "2023-09-11 Monday"

Replace-Text "70÷2=" "33÷3="
Replace-Text "10÷4=" "32÷4="
Replace-Text "96÷5=" "17÷5="
Replace-Text "32÷7=" "31÷9="
Replace-Text "32÷8=" "90÷2="

Replace-Text "93÷2=" "85÷4="
Replace-Text "73÷8=" "31÷5="
Replace-Text "30÷3=" "61÷2="
Replace-Text "47÷5=" "26÷9="
Replace-Text "15÷8=" "98÷9="

Replace-Text "39÷9=" "99÷4="
Replace-Text "28÷9=" "27÷5="
Replace-Text "88÷9=" "97÷8="
Replace-Text "87÷3=" "16÷7="
Replace-Text "15÷3=" "46÷2="

Replace-Text "18÷8=" "22÷3="
Replace-Text "16÷4=" "84÷2="
Replace-Text "34÷7=" "62÷3="
Replace-Text "41÷9=" "35÷3="
Replace-Text "74÷7=" "43÷3="

Replace-Text "48÷8=" "70÷2="
Replace-Text "31÷4=" "49÷3="
Replace-Text "20÷5=" "56÷5="
Replace-Text "19÷8=" "85÷7="
Replace-Text "91÷7=" "30÷4="
